$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16:21 down to 17:22
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new weekly record (same fixed fields as the
# surrounding Espinaca / Terminal Hortofrutícola Agro Chillán rows).
$row = 16

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2022 -Month 10 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112012
$ws.Cells.Item($row, 7).Value = "Espinaca"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 6500
$ws.Cells.Item($row, 12).Value = 7000
$ws.Cells.Item($row, 13).Value = 6750
$ws.Cells.Item($row, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value = 675
$ws.Cells.Item($row, 17).Value = 10
$ws.Cells.Item($row, 18).Value = "Hortaliza"
